$d = $word.ActiveDocument

# --- Change 1: header date range (title paragraph) ---
$d.Content.Find.Execute(
    "Du 01/09/2022 au 26/04/2023", $true, $false, $false, $false, $false,
    $true, 1, $false, "Du 9/29/2022 au 10/1/2022", 2) | Out-Null

# --- Change 2: "Gains/pertes globals" summary sentence ---
$d.Content.Find.Execute(
    "Du 2022-09-01T00:00:00.000Z au 2023-04-26T00:00:00.000Z, vous avez gagné 0euro.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Du 2022-09-29T00:00:00.000Z au 2022-10-01T00:00:00.000Z, vous avez gagné 0euro.", 2) | Out-Null

# --- Change 3: rebuild the "Investissements" table cell (SOL row) ---
# The cell currently holds 6 paragraphs (one combined BUY/BUY/SELL paragraph
# followed by five repeated SELL paragraphs). It must become exactly two
# paragraphs: one BUY line and one SELL line.
#
# NOTE: do not call $d.Tables.Item(...) before walking $d.Paragraphs in this
# host -- doing so corrupts the live paragraph iterator. Locate the target
# paragraphs purely via the document-level Paragraphs collection instead.
$solIndex = 0
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text.StartsWith("SOL")) {
        $solIndex = $idx
    }
}

$firstIndex = $solIndex + 1
$k = $firstIndex
while ($d.Paragraphs.Item($k + 1).Range.Text.StartsWith("SELL")) {
    $k = $k + 1
}
$lastIndex = $k

# Delete every paragraph in the cell after the first one.
for ($j = $lastIndex; $j -gt $firstIndex; $j--) {
    $d.Paragraphs.Item($j).Range.Delete() | Out-Null
}

# Rewrite the surviving (first) paragraph with the BUY line.
$firstPara = $d.Paragraphs.Item($firstIndex)
$firstPara.Range.Text = "BUY -> 10 sol à 10 eur `n"

# Append a brand-new, cleanly-formatted paragraph for the SELL line.
$endRange = $firstPara.Range
$endRange.Collapse(0) | Out-Null
$endRange.InsertAfter([char]13) | Out-Null

$secondPara = $d.Paragraphs.Item($firstIndex + 1)
$secondPara.Range.Text = "SELL -> 10 sol pour 20 € `n"
